# Replace the "books" sample sheet with a small "Location / Population" table
# (ExcelObject sample data), matching the new sampledata.xlsx contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - remove the old title/author/pages book list entirely.
$ws.Cells.Clear()

# Headers
$ws.Range("A1").Value = "Location"
$ws.Range("B1").Value = "Population"

# Rows
$ws.Range("A2").Value = "Naboo"
$ws.Range("B2").Value = 1000000

$ws.Range("A3").Value = "Alderaan"
$ws.Range("B3").Value = 0

$ws.Range("A4").Value = "Tatooine"
$ws.Range("B4").Value = 5000000

$ws.Range("A5").Value = "Coruscant"
$ws.Range("B5").Value = 10000000000

# Population numbers (except the 0 row) use a thousands-separated integer format
$ws.Range("B2").NumberFormat = "#,##0"
$ws.Range("B4:B5").NumberFormat = "#,##0"

# Column widths set to match the saved workbook
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(2).ColumnWidth = 17.17
$ws.Columns.Item(3).ColumnWidth = 17

# Leave selection where Excel left it when the file was saved
[void]$ws.Range("B6").Select()
